# Completes the comment of the PEXTRA table 126 variables in the
# (pre) identified missing file (issue #450).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Simple text appends: add "part of MFPPHY" / "part of MFP3D" markers
#    to the existing PEXTRA table 126 comments.
# ---------------------------------------------------------------------

$ws.Range("H37").Value = "COSP grib 126.40   CVEXTR2(1)='CALIPSO_LCC'                part of MFPPHY"
$ws.Range("H38").Value = "COSP grib 126.41   CVEXTR2(2)='CALIPSO_MCC'               part of MFPPHY"
$ws.Range("H39").Value = "COSP grib 126.42   CVEXTR2(3)='CALIPSO_HCC'                part of MFPPHY"
$ws.Range("H40").Value = "COSP grib 126.43   CVEXTR2(4)='CALIPSO_TCC'                part of MFPPHY"
$ws.Range("H41").Value = "COSP grib 126.44   CVEXTR2(5)='ISCCP_TOTALCLD'       part of MFPPHY"
$ws.Range("H42").Value = "COSP grib 126.45   CVEXTR2(6)='ISCCP_MEANPTOP'           part of MFPPHY"
$ws.Range("H43").Value = "COSP grib 126.46   CVEXTR2(7)='ISCCP_MEANALBCLD'       part of MFPPHY"

$ws.Range("H46").Value = 'grib 126.73                          part of MFPPHY    Available from double radiation call in IFS. PEXTRA issue #403   aerosol free'
$ws.Range("H47").Value = 'grib 126.72                          part of MFPPHY    Available from double radiation call in IFS. PEXTRA issue #403   aerosol free'
$ws.Range("H48").Value = 'grib 128.212-126.069     part of MFPPHY    Available from double radiation call in IFS. PEXTRA issue #403   aerosol free'

$ws.Range("H51").Value = "Grib 126.20 / 126.22        part of MFP3D        In namelist.ifs.cloudact+diag.sh  CVEXTRA(1)='CDNC' which is a PEXTRA variable."

$ws.Range("H54").Value = 'grib 126.95                                                                      part of MFP3D        Available in IFS: T-tendency from radiation: grib 128.95'
$ws.Range("H55").Value = 'grib 126.105                                                                   part of MFP3D        Available in IFS: T-tendency from convection : grib 128.105'
$ws.Range("H56").Value = 'Grib 126.94 + 126.99 + 126.106 + 126.110       part of MFP3D        Adding all the q-tendencies, thus: grib 128.94 + 128.99 + 128.106 + 128.110.  Alternatively, in IFS: just estimating the delta q per month. So far no direct grib code for the totoal q-tendency found'
$ws.Range("H57").Value = 'grib 126.106                                                                   part of MFP3D        Available in IFS: q-tendency from convection: grib 128.106'
$ws.Range("H58").Value = 'grib 126.99 + 126.106 + 126.110                           part of MFP3D        Adding all the q-tendencies without advection, thus: grib 128.99 + 128.106 + 128.110.'

# ---------------------------------------------------------------------
# 2. Rich-text comments for sftgrf (H61), sncIs (H62), sftgif (H63) and
#    areacellg (H64): prefix "To be implemented:" plus a red grib code,
#    followed by "part of MFPPHY" and the original comment text.
# ---------------------------------------------------------------------

function Set-RichComment($cell, $runs) {
    $full = ""
    foreach ($run in $runs) { $full += $run.Text }
    $cell.Value = $full
    $pos = 1
    foreach ($run in $runs) {
        $len = $run.Text.Length
        if ($run.Color -ne $null) {
            $cell.Characters($pos, $len).Font.Color = $run.Color
        }
        $pos += $len
    }
}

$newRed  = 1972430   # RGB(0xCE,0x18,0x1E) -> FFCE181E
$origRed = 255        # RGB(0xFF,0x00,0x00) -> FFFF0000

Set-RichComment $ws.Range("H61") @(
    @{ Text = 'To be implemented: '; Color = $null },
    @{ Text = ' grib 126.30'; Color = $newRed },
    @{ Text = '  part of MFPPHY   For Greenland this is the same as above sftgif. We do not have Antarctic ice sheet.'; Color = $null }
)

Set-RichComment $ws.Range("H62") @(
    @{ Text = 'To be implemented:  '; Color = $null },
    @{ Text = 'grib 126.31'; Color = $newRed },
    @{ Text = '  part of MFPPHY   Not available in IFS. Although it could be calculated from tile fractions and written out as extra output'; Color = $null }
)

Set-RichComment $ws.Range("H63") @(
    @{ Text = 'To be implemented:  '; Color = $null },
    @{ Text = 'grib 126.32'; Color = $newRed },
    @{ Text = '  part of MFPPHY   This is the land ice mask and will be an extra variable in IFS ('; Color = $null },
    @{ Text = 'thomas: via PEXTRA?'; Color = $origRed },
    @{ Text = ')'; Color = $null }
)

Set-RichComment $ws.Range("H64") @(
    @{ Text = 'To be implemented:  '; Color = $null },
    @{ Text = 'grib 126.34'; Color = $newRed },
    @{ Text = '  part of MFPPHY   Available in PISM. This is the ice sheet mask (in fraction) defined in the ice sheet model grid'; Color = $null }
)

# ---------------------------------------------------------------------
# 3. Row height of the rows that now contain the richer, taller comments.
# ---------------------------------------------------------------------

$ws.Rows.Item(61).RowHeight = 14.25
$ws.Rows.Item(62).RowHeight = 14.25
$ws.Rows.Item(64).RowHeight = 14.25

# ---------------------------------------------------------------------
# 4. Restore the view/selection state recorded in the workbook.
# ---------------------------------------------------------------------

$ws.Range("E57").Select()
